$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, "wenden", "none", "none")
    ,@(3, "sehen", "dog/dog020.jpg", "dog")
    ,@(4, "tauschen", "flower/flower030.jpg", "flower")
    ,@(5, "rufen", "none", "none")
    ,@(6, "fügen", "flower/flower019.jpg", "flower")
    ,@(7, "altern", "dog/dog003.jpg", "dog")
    ,@(8, "enden", "none", "none")
    ,@(9, "triefen", "dog/dog012.jpg", "dog")
    ,@(10, "zögern", "dog/dog009.jpg", "dog")
    ,@(11, "opfern", "none", "none")
    ,@(12, "kümmern", "flower/flower032.jpg", "flower")
    ,@(13, "ehren", "dog/dog027.jpg", "dog")
    ,@(14, "schalten", "none", "none")
    ,@(15, "treiben", "dog/dog014.jpg", "dog")
    ,@(16, "tropfen", "flower/flower011.jpg", "flower")
    ,@(17, "klagen", "none", "none")
    ,@(18, "kehren", "flower/flower005.jpg", "flower")
    ,@(19, "bauen", "dog/dog005.jpg", "dog")
    ,@(20, "sparen", "none", "none")
    ,@(21, "wundern", "flower/flower026.jpg", "flower")
    ,@(22, "wachsen", "flower/flower029.jpg", "flower")
    ,@(23, "hören", "none", "none")
    ,@(24, "schleppen", "dog/dog025.jpg", "dog")
    ,@(25, "fragen", "flower/flower031.jpg", "flower")
    ,@(26, "drohen", "none", "none")
    ,@(27, "buchen", "dog/dog006.jpg", "dog")
    ,@(28, "packen", "flower/flower012.jpg", "flower")
    ,@(29, "orten", "none", "none")
    ,@(30, "klingen", "flower/flower013.jpg", "flower")
    ,@(31, "hoffen", "flower/flower010.jpg", "flower")
    ,@(32, "weigern", "none", "none")
    ,@(33, "schneiden", "dog/dog017.jpg", "dog")
    ,@(34, "dürfen", "dog/dog026.jpg", "dog")
    ,@(35, "ächzen", "none", "none")
    ,@(36, "pflegen", "flower/flower002.jpg", "flower")
    ,@(37, "platzen", "dog/dog010.jpg", "dog")
    ,@(38, "kosten", "none", "none")
    ,@(39, "knien", "flower/flower022.jpg", "flower")
    ,@(40, "töten", "flower/flower027.jpg", "flower")
    ,@(41, "dauern", "none", "none")
    ,@(42, "husten", "flower/flower018.jpg", "flower")
    ,@(43, "rasen", "dog/dog001.jpg", "dog")
    ,@(44, "stören", "none", "none")
    ,@(45, "machen", "dog/dog030.jpg", "dog")
    ,@(46, "mühen", "dog/dog004.jpg", "dog")
    ,@(47, "bremsen", "none", "none")
    ,@(48, "beißen", "dog/dog023.jpg", "dog")
    ,@(49, "ärgern", "flower/flower003.jpg", "flower")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
